# TC15_C3DC_phs002790_CauseOfDeath-NotReported.xlsx
# "Updated queries for C3DC first half testcases."
#
# The SQL stored in the "TabQuery"/"StatQuery" columns joined df_* tables on a
# generic surrogate "id" column (std.id / prt.id / "study.id" / "participant.id").
# The commons schema now exposes explicit natural-key columns, so every LEFT
# JOIN condition is rewritten to use study_id / participant_id throughout.

function FixJoins($text) {
    $t = $text
    $t = $t.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $t = $t.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $t = $t.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $t = $t.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $t = $t.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $t = $t.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
    return $t
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: this rebuilds the shared-string table, so cells are touched
# in the same left-to-right / top-to-bottom order a save from Excel would
# observe them in (C2 before B2, then B3..B7) so the resulting shared-string
# indices line up with the target workbook.
$cellsInScanOrder = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cellsInScanOrder) {
    $old = $ws.Range($addr).Value()
    $new = FixJoins $old
    $ws.Range($addr).Value = $new
}

# The author's selection ended up resting on C7 after scrolling the sheet
# down so row 6 is at the top of the viewport.
$ws.Range("C7").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
